$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (34-36) following the same pattern as the
# existing data rows: regcntr_id=10005, incrementing usr_id values,
# and the same lang_code/is_active/cr_by/cr_dtimes/eff_dtimes values.
$newRows = @(
    @(10005, 110033),
    @(10005, 110034),
    @(10005, 110035)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update the selection to mirror the post-edit state captured in the diff
$ws.Range("A37:XFD1048576").Select()
